$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-08-05 Monday" "2024-08-06 Tuesday"

Replace-Text "394×9=" "276×7="
Replace-Text "791×9=" "622×3="
Replace-Text "195×9=" "656×5="
Replace-Text "693×6=" "413×8="
Replace-Text "432×4=" "324×8="
Replace-Text "565×7=" "543×2="
Replace-Text "353×7=" "733×8="
Replace-Text "862×5=" "830×5="
Replace-Text "813×6=" "216×7="
Replace-Text "217×9=" "724×9="
Replace-Text "722×7=" "422×5="
Replace-Text "962×3=" "447×8="
Replace-Text "320×9=" "175×3="
Replace-Text "490×5=" "272×8="
Replace-Text "235×6=" "363×4="
Replace-Text "207×2=" "616×3="
Replace-Text "479×8=" "692×3="
Replace-Text "154×9=" "335×3="
Replace-Text "941×9=" "386×7="
Replace-Text "520×6=" "643×3="
Replace-Text "449×7=" "741×9="
Replace-Text "919×6=" "990×3="
Replace-Text "947×4=" "663×6="
Replace-Text "450×6=" "507×4="
Replace-Text "994×6=" "714×7="
